$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.017.10"
$ws.Range("E2").Value = "  +2.08%  "

$ws.Range("D3").Value = "1.909.79"
$ws.Range("E3").Value = "  +2.37%  "

$ws.Range("E4").Value = "  -0.78%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.33"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  +1.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07362"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9338"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.82"
$ws.Range("E11").Value = "  +0.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07772"
$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").Value = "1.898.62"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.499"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.627"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.91"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"

$ws.Range("D20").Value = "28.036.66"
$ws.Range("E20").Value = "  +1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.78"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.167"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").Value = "2.144.79"
$ws.Range("E23").Value = "  +1.15%  "

$ws.Range("E24").Value = "  +1.77%  "

$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.915"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.141"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.01"
$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.961"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08944"
$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.296"
$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.264"
$ws.Range("E33").Value = "  +4.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7780"
$ws.Range("E34").Value = "  +3.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.685"
$ws.Range("E35").Value = "  +1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.632"
$ws.Range("E36").Value = "  -3.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02056"
$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.111"
$ws.Range("E38").Value = "  -0.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05313"
$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.010"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5495"
$ws.Range("E41").Value = "  +3.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.021"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.484"
$ws.Range("E44").Value = "  -0.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").Value = "  +0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4829"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.90"
$ws.Range("E47").Value = "  +4.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.92"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06072"
$ws.Range("E51").Value = "  -0.11%  "
